# Remove the RG (ID card) reference from the certificate text and
# tidy up the "São Carlos, {{DATA}}" run split, per the commit
# "remove RG to user and models certificate".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)          # "Rectangle 5" - body paragraph shape
$tf  = $shp.TextFrame
$tr  = $tf.TextRange

$masc = [char]0x00BA              # 'º' (masculine ordinal indicator)

# The shape has <a:spAutoFit/> (resize-shape-to-fit-text), so trimming
# text below would otherwise shrink its saved height. Remember the
# original size so it can be restored once all edits are done.
$origHeight = $shp.Height
$origWidth  = $shp.Width

# ---------------------------------------------------------------
# 1) "...portador do RG nº {{RG}} e CPF nº {{CPF}}, ..."
#        -> "...portador do CPF nº {{CPF}}, ..."
# ---------------------------------------------------------------

# 1a. Drop the "{{RG}} " run entirely.
$full = $tr.Text
$needle = "{{RG}} "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = ""
}

# 1b. "e CPF nº " -> "nº " (drop the leading "e CPF ").
$full = $tr.Text
$needle = "e CPF n" + $masc + " "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = "n" + $masc + " "
}

# 1c. Strip the now-redundant "RG nº " out of "portador do RG nº ",
#     leaving "portador do ".
$full = $tr.Text
$needle = "RG n" + $masc + " "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = ""
}

# 1d. Split "portador do " into separate "portador " / "do " runs.
$full = $tr.Text
$needle = "do "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = "do "
}

# 1e. Insert "CPF " right after "do " so the text reads
#     "portador do CPF nº {{CPF}}, ".
$full = $tr.Text
$needle = "do "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    [void]$rng.InsertAfter("CPF ")
}

# 1f. Force "CPF " onto its own run (separate from "nº ").
$full = $tr.Text
$needle = "CPF "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = "CPF "
}

# ---------------------------------------------------------------
# 2) "São Carlos" + ", " -> single run "São Carlos, "
# ---------------------------------------------------------------
$full = $tr.Text
$needle = "S" + [char]0x00E3 + "o Carlos, "
$idx = $full.IndexOf($needle)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $needle.Length)
    $rng.Text = $needle
}

# Restore the shape's original autofit-computed size so that trimming
# the text doesn't also shrink the box in the saved file.
$shp.Height = $origHeight
$shp.Width  = $origWidth

Write-Host "Final text:"
Write-Host $tr.Text
